$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ------------------------------------------------------
# New columns F (RampTime(Hrs)) and G (Limits) are appended, and C1's
# header text changes from RoundTripEfficiency to AvgEfficiency.
$ws.Range("G1").Value = "Limits"
$ws.Range("F1").Value = "RampTime(Hrs)"
$ws.Range("C1").Value = "AvgEfficiency"

# --- Row 2 (Battery) ---------------------------------------------------
$ws.Range("D2").Value = 0.8
$ws.Range("E2").Value = "inf"
$ws.Range("F2").Value = 0.25
$ws.Range("G2").Value = 0.7

# --- Row 3 (Hydrogen) ---------------------------------------------------
$ws.Range("E3").Value = "inf"
$ws.Range("F3").Value = 0.25
$ws.Range("G3").Value = 0.8

# --- Row 4 (Hydro) ---------------------------------------------------
$ws.Range("B4").Value = 5000
$ws.Range("D4").Value = 0.0001
$ws.Range("F4").Value = 6
$ws.Range("G4").Value = 0.5

# --- Column widths (bestFit, matches Excel's auto-fit after edits) -----
$ws.Columns.Item(3).ColumnWidth = 11
$ws.Columns.Item(4).ColumnWidth = 6.333333333333334
$ws.Columns.Item(5).ColumnWidth = 10.666666666666666
$ws.Columns.Item(6).ColumnWidth = 12.833333333333332
$ws.Columns.Item(7).ColumnWidth = 5

# --- Selection -----------------------------------------------------
$null = $ws.Range("C2").Select()
